$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: date "January 20, 2018" -> "January 11, 2018"
# The original run layout splits the date oddly across 4 runs
# ("January" / " 20, 201" / "8" / " "). The target layout instead splits
# it as "January " / "11" / ", 2018 " (3 runs, all w:rPr/ empty).
# We replace just the "20" characters, then nudge formatting on that
# sub-range (set+clear Bold) so the engine keeps it as its own run
# (with an explicit empty <w:rPr/>) instead of silently re-merging it
# into its neighbours.
# ---------------------------------------------------------------------

$dateFind = $d.Content
$found = $dateFind.Find.Execute("January 20, 2018", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "", 0)
$dateStart = $dateFind.Start

# "January " occupies the first 8 chars, "20" the next 2, ", 2018 " after that.
$numRange = $d.Range($dateStart + 8, $dateStart + 10)
$numRange.Text = "11"

# Force the edited run to keep its own identity / explicit empty rPr.
$numRange = $d.Range($dateStart + 8, $dateStart + 10)
$numRange.Font.Bold = 1
$numRange.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 2: collapse the separate address runs into single runs (no
# visible text changes) - "Suite 160 / Box 650 / Seattle..." become one
# run, and "Current Physical Address: / 8244 .../ " become one run,
# joined by the existing <w:br/>.
# ---------------------------------------------------------------------

$addrFind = $d.Content
$found = $addrFind.Find.Execute("Suite 160", $false, $false, $false, $false, $false,
                                 $true, 1, $false, "", 0)
$suiteStart = $addrFind.Start

$addrFind2 = $d.Content
$found2 = $addrFind2.Find.Execute("98033", $false, $false, $false, $false, $false,
                                   $true, 1, $false, "", 0)
$zipEnd = $addrFind2.End + 1   # include the trailing space after 98033

# Touch the span starting one character inside the "Suite 160" run (not
# exactly at its shared boundary with the preceding "  .  " run) so the
# auto-merge pass only coalesces the runs we want merged, leaving the
# "  .  " run (and the bookmark around it) untouched.
$mergeRange = $d.Range($suiteStart + 1, $zipEnd)
$original = $mergeRange.Text
$mergeRange.Text = $original + "X"
$mergeRange2 = $d.Range($suiteStart + 1, $suiteStart + 1 + $original.Length + 1)
$mergeRange2.Text = $original

# Restore explicit empty <w:rPr/> on the merged run and on the
# preceding "  .  " run (both lost their rPr element when touched).
$preRange = $d.Range($suiteStart - 5, $suiteStart)
$preRange.Font.Bold = 1
$preRange.Font.Bold = 0

$bigRange = $d.Range($suiteStart, $zipEnd)
$bigRange.Font.Bold = 1
$bigRange.Font.Bold = 0
